$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new column R values for year 2022
$ws.Range("R4").Value = 2022
$ws.Range("R5").Value = 8.6821914120339212
$ws.Range("R6").Value = 12.221423436376707

# Copy styles from the Q column (previous last column) to the new R column
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122)

$ws.Range("Q6").Copy()
$ws.Range("R6").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Update selection to match target active cell
$ws.Range("S4").Select()
